$d = $word.ActiveDocument

# The sample JCL snippet showed a misleading "DATA(16305,16305)" value.
# Revise it to "DATE(2023305,2023305)" so readers aren't misled.
$d.Content.Find.Execute("DATA(16305,16305)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DATE(2023305,2023305)", 2) | Out-Null
